# "fixed perubahan besar pada nomenklatur"
# Adds Kode Urusan / Kode Program / Kode Kegiatan columns ahead of the
# existing Kode Sub Kegiatan / Sub Kegiatan / Tahun Perubahan header, drops
# the old sample data row, and leaves two formatted-but-empty input cells
# behind for the (now shifted) Kode Sub Kegiatan / Sub Kegiatan columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Grab the data-row style (s="1", small grey Arial) that used to live on
#    B2/C2 and stamp it onto the new E2:F2 placeholder cells *before* the
#    old sample row gets wiped out below.
$ws.Range("C2").Copy()
$ws.Range("E2:F2").PasteSpecial(-4122)

# 2. Shift the existing header labels from B1:D1 right, into E1:G1.
$ws.Range("G1").Value = $ws.Range("D1").Formula
$ws.Range("F1").Value = $ws.Range("C1").Formula
$ws.Range("E1").Value = $ws.Range("B1").Formula

# 3. Drop the old sample/demo data row (A2:C2).
$ws.Range("A2:C2").Clear()

# 4. Populate the three newly inserted header columns.
$ws.Range("B1").Value = "Kode Urusan"
$ws.Range("C1").Value = "Kode Program"
$ws.Range("D1").Value = "Kode Kegiatan"

# 5. Re-apply the column widths for the (now shifted) columns.
$ws.Columns("B").ColumnWidth = 10.5
$ws.Columns("C").ColumnWidth = 11.666666666666666
$ws.Columns("D").ColumnWidth = 11.666666666666666
$ws.Columns("E").ColumnWidth = 15.166666666666666
$ws.Columns("F").ColumnWidth = 45.33333333333333
$ws.Columns("G").ColumnWidth = 14.333333333333332

# 6. Match the saved cursor/selection position.
$ws.Range("J6").Select()
